$d = $word.ActiveDocument

function Get-ParagraphIndexContainingPosition($doc, $pos) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $pos -and $p.Range.End -gt $pos) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "Personale" section: fill in the empty paragraph right after the
#    "Personale" heading with the staff names, then anchor a "_GoBack"
#    bookmark right after the final word ("Vitale"), mirroring the target
#    markup:
#      <w:r><w:t>...Emilio Meroni e Benedetta Vitale</w:t></w:r>
#      <w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
# ---------------------------------------------------------------------------

$headingScope = $d.Content
$headingScope.Find.Execute("Personale", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$personaleIdx = Get-ParagraphIndexContainingPosition $d $headingScope.Start
$personaleEmptyPara = $d.Paragraphs.Item($personaleIdx + 1)

$insertRange = $d.Range($personaleEmptyPara.Range.Start, $personaleEmptyPara.Range.Start)
# Trailing "~" is a throw-away placeholder character. Adding a collapsed
# bookmark exactly at a run/paragraph boundary is unreliable, so the
# placeholder gives a safe "mid run" insertion point for the bookmark; it
# is deleted again right after.
$insertRange.InsertAfter("Il personale che ha lavorato a questo programma sono: Emilio Meroni e Benedetta Vitale~")

$placeholderPos = $insertRange.End - 1
$bookmarkRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Text = ""

# ---------------------------------------------------------------------------
# 2) "garanzia e qualità" section: fill the empty paragraph after the
#    heading with a short summary sentence, then add a second, new
#    paragraph with the longer explanation right after it.
# ---------------------------------------------------------------------------

$qualitaScope = $d.Content
$qualitaScope.Find.Execute("garanzia e qualit" + [char]0x00E0, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$qualitaIdx = Get-ParagraphIndexContainingPosition $d $qualitaScope.Start
$qualitaEmptyPara = $d.Paragraphs.Item($qualitaIdx + 1)

$qualitaInsertRange = $d.Range($qualitaEmptyPara.Range.Start, $qualitaEmptyPara.Range.Start)

$shortText = 'Il programma ' + [char]0x00E8 + ': semplice, intuibile, veloce e affidabile.'

$longText = [char]0x00C8 + ' semplice e intuibile perch' + [char]0x00E9 + ' ha poche scritte di facile comprensione, anche per chi si approccia al programma per la prima volta al programma; ci sono poche schermate che hanno al loro interno tutto il necessario per le operazioni che si desiderano; inoltre i bottoni e i colori rendono il programma semplice e intuibile. ' + [char]0x00C8 + ' veloce in quanto il programma ha poche sezioni in cui bisogna utilizzare la testiera, mentre il programma ' + [char]0x00E8 + ' pi' + [char]0x00F9 + ' focalizzato sull' + [char]0x2019 + 'uso dei bottoni che lo rendono molto pi' + [char]0x00F9 + ' veloce nelle azioni che si dovranno svolgere. Il programma ' + [char]0x00E8 + ' affidabile dato che, con le schermate pop-up di conferma, c' + [char]0x2019 + [char]0x00E8 + ' una riduzione dei possibili errori che si possono creare nell' + [char]0x2019 + 'uso del programma. '

$qualitaInsertRange.InsertAfter($shortText + [char]13 + $longText)

Write-Output "done"
